$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "001" to "004"; K2 stays "001" (unchanged).
# Leading apostrophe forces text storage (avoids Excel's auto numeric
# coercion of "004"); resetting the style afterwards avoids leaving a
# stray quote-prefix style applied to the cell.
$ws.Range("J2").Value = "'004"
$ws.Range("J2").Style = "Normal"

# N2: report date changes
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Numeric columns updated with new figures
$ws.Range("O2").Value = 1839759362.45
$ws.Range("P2").Value = 599377150.52
$ws.Range("Q2").Value = 368081553.34
$ws.Range("R2").ClearContents()
$ws.Range("S2").Value = 60629975.98
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 84337111.62
$ws.Range("V2").ClearContents()
$ws.Range("W2").Value = 1079788221.18
$ws.Range("X2").Value = 176490984.3
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").Value = 441938.75
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").Value = 759971141.27
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").Value = 77.9801347275
$ws.Range("AG2").Value = 58.6918182464
